$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "'320018539885"
$ws.Range("P3").Value = "'320018539896"
$ws.Range("P4").Value = "'320018539922"
$ws.Range("P5").Value = "'320018539933"
$ws.Range("Q4").Value = "'$49.70"
$ws.Range("R4").Value = "PASS"

$ws.Range("P2").Style = "Normal"
$ws.Range("P3").Style = "Normal"
$ws.Range("P4").Style = "Normal"
$ws.Range("P5").Style = "Normal"
$ws.Range("Q4").Style = "Normal"
